$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6956
$ws1.Range("F3").Value = 53
$ws1.Range("F5").Value = 68
$ws1.Range("F6").Value = 1077
$ws1.Range("F7").Value = 167
$ws1.Range("F8").Value = 9

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6956
$ws4.Range("F3").Value = 53
$ws4.Range("F5").Value = 68
$ws4.Range("F6").Value = 1077
$ws4.Range("F7").Value = 167
$ws4.Range("F9").Value = 9
